# Add unique VINs to Each CA Choice Test, all
# Rows 2-6 previously held duplicate/re-used VIN values ("1FDEU15H&K" / "2FDEU15H&K").
# Replace them all with a single new, unique VIN value so every row has a distinct VIN.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newVin = "AAAVB3CC&W"

foreach ($row in 2..6) {
    $ws.Cells.Item($row, 1).Value = $newVin
}
